# Regenerate merged AHB files: rename header columns, freeze header row,
# and wrap the data range in a table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AHB-Diff")

# 1) Rename the "_old"/"_new" header labels to "_FV2404"/"_FV2410"
$headers = @(
    "Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404",
    "Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410",
    "Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2) Freeze the header row (split after row 1)
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# 3) Wrap A1:U86 in a table, preserving the header row's existing formatting
#    (stash it via copy/paste so the new table doesn't bake a header dxf).
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A200:U200")
$headerRange.Copy()
$scratch.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$headerRange.ClearFormats()

$dataRange = $ws.Range("A1:U86")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$scratch.Copy()
$headerRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$scratch.ClearFormats()
$scratch.ClearContents()

$ws.Range("A1").Select() | Out-Null
